$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("热源")

# --- Insert the two new rows first (keeps formatting/height consistent with
#     the rows they sit next to), then fill in values afterwards so the new
#     shared-string entries get allocated in the same order as upstream
#     ("含水坩埚" first -> index 90, then "含岩浆坩埚" -> index 91). ---

# Row for "含岩浆坩埚" (1400) goes above the existing "水" row (row 6).
$ws.Rows("6:6").Insert()
$ws.Range("A7:B7").Copy()
$ws.Range("A6:B6").PasteSpecial(-4122)
$ws.Rows("6:6").RowHeight = $ws.Rows("7:7").RowHeight

# Row for "含水坩埚" (15) goes below "水" (now row 7), above "冰" (now row 8).
$ws.Rows("8:8").Insert()
$ws.Range("A9:B9").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)
$ws.Rows("8:8").RowHeight = $ws.Rows("9:9").RowHeight

# Fill in values in shared-string allocation order: 含水坩埚 (90) then 含岩浆坩埚 (91).
$ws.Range("A8").Value = "含水坩埚"
$ws.Range("B8").Value = 15
$ws.Range("A6").Value = "含岩浆坩埚"
$ws.Range("B6").Value = 1400

$excel.CutCopyMode = $false

# --- Update the active sheet / selection to match the new view state ---
$ws.Activate()
$ws.Range("C6").Select()
